$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Clcf1"
$ws.Cells.Item(2,3).Value = "Crlf1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 1.798783666666667
$ws.Cells.Item(2,8).Value = 5.396351
$ws.Cells.Item(2,9).Value = 0.2319744053785674
$ws.Cells.Item(2,10).Value = 0.2319744053785674
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.082925
$ws.Cells.Item(2,14).Value = 0.248775
$ws.Cells.Item(2,15).Value = 0.003066867285585202
$ws.Cells.Item(2,16).Value = 0.003066867285585203
$ws.Cells.Item(2,17).Value = 0.1491641355583333
$ws.Cells.Item(2,18).Value = 1.342477220025
$ws.Cells.Item(2,19).Value = 0.0007114347149486084
$ws.Cells.Item(2,20).Value = 0.0007114347149486084

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Clcf1"
$ws.Cells.Item(3,3).Value = "Crlf1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 1.798783666666667
$ws.Cells.Item(3,8).Value = 5.396351
$ws.Cells.Item(3,9).Value = 0.2319744053785674
$ws.Cells.Item(3,10).Value = 0.2319744053785674
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 17.662076
$ws.Cells.Item(3,14).Value = 52.986228
$ws.Cells.Item(3,15).Value = 0.653207634367435
$ws.Cells.Item(3,16).Value = 0.653207634367435
$ws.Cells.Item(3,17).Value = 31.77025382822533
$ws.Cells.Item(3,18).Value = 285.932284454028
$ws.Cells.Item(3,19).Value = 0.1515274525711264
$ws.Cells.Item(3,20).Value = 0.1515274525711264

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Clcf1"
$ws.Cells.Item(4,3).Value = "Crlf1"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 1.798783666666667
$ws.Cells.Item(4,8).Value = 5.396351
$ws.Cells.Item(4,9).Value = 0.2319744053785674
$ws.Cells.Item(4,10).Value = 0.2319744053785674
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 9.293991
$ws.Cells.Item(4,14).Value = 27.881973
$ws.Cells.Item(4,15).Value = 0.3437254983469798
$ws.Cells.Item(4,16).Value = 0.3437254983469798
$ws.Cells.Item(4,17).Value = 16.717879208947
$ws.Cells.Item(4,18).Value = 150.460912880523
$ws.Cells.Item(4,19).Value = 0.0797355180924924
$ws.Cells.Item(4,20).Value = 0.0797355180924924

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Clcf1"
$ws.Cells.Item(5,3).Value = "Crlf1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.782700333333333
$ws.Cells.Item(5,8).Value = 5.348101
$ws.Cells.Item(5,9).Value = 0.2299002695301921
$ws.Cells.Item(5,10).Value = 0.2299002695301921
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.082925
$ws.Cells.Item(5,14).Value = 0.248775
$ws.Cells.Item(5,15).Value = 0.003066867285585202
$ws.Cells.Item(5,16).Value = 0.003066867285585203
$ws.Cells.Item(5,17).Value = 0.1478304251416667
$ws.Cells.Item(5,18).Value = 1.330473826275
$ws.Cells.Item(5,19).Value = 0.0007050736155693666
$ws.Cells.Item(5,20).Value = 0.0007050736155693666

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Clcf1"
$ws.Cells.Item(6,3).Value = "Crlf1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.782700333333333
$ws.Cells.Item(6,8).Value = 5.348101
$ws.Cells.Item(6,9).Value = 0.2299002695301921
$ws.Cells.Item(6,10).Value = 0.2299002695301921
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 17.662076
$ws.Cells.Item(6,14).Value = 52.986228
$ws.Cells.Item(6,15).Value = 0.653207634367435
$ws.Cells.Item(6,16).Value = 0.653207634367435
$ws.Cells.Item(6,17).Value = 31.48618877255867
$ws.Cells.Item(6,18).Value = 283.3756989530279
$ws.Cells.Item(6,19).Value = 0.1501726112002525
$ws.Cells.Item(6,20).Value = 0.1501726112002525

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Clcf1"
$ws.Cells.Item(7,3).Value = "Crlf1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.782700333333333
$ws.Cells.Item(7,8).Value = 5.348101
$ws.Cells.Item(7,9).Value = 0.2299002695301921
$ws.Cells.Item(7,10).Value = 0.2299002695301921
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 9.293991
$ws.Cells.Item(7,14).Value = 27.881973
$ws.Cells.Item(7,15).Value = 0.3437254983469798
$ws.Cells.Item(7,16).Value = 0.3437254983469798
$ws.Cells.Item(7,17).Value = 16.568400853697
$ws.Cells.Item(7,18).Value = 149.115607683273
$ws.Cells.Item(7,19).Value = 0.07902258471437025
$ws.Cells.Item(7,20).Value = 0.07902258471437025

$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Clcf1"
$ws.Cells.Item(8,3).Value = "Crlf1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 4.172749333333334
$ws.Cells.Item(8,8).Value = 12.518248
$ws.Cells.Item(8,9).Value = 0.5381253250912406
$ws.Cells.Item(8,10).Value = 0.5381253250912404
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.082925
$ws.Cells.Item(8,14).Value = 0.248775
$ws.Cells.Item(8,15).Value = 0.003066867285585202
$ws.Cells.Item(8,16).Value = 0.003066867285585203
$ws.Cells.Item(8,17).Value = 0.3460252384666667
$ws.Cells.Item(8,18).Value = 3.1142271462
$ws.Cells.Item(8,19).Value = 0.001650358955067228
$ws.Cells.Item(8,20).Value = 0.001650358955067227

$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Clcf1"
$ws.Cells.Item(9,3).Value = "Crlf1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 4.172749333333334
$ws.Cells.Item(9,8).Value = 12.518248
$ws.Cells.Item(9,9).Value = 0.5381253250912406
$ws.Cells.Item(9,10).Value = 0.5381253250912404
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 17.662076
$ws.Cells.Item(9,14).Value = 52.986228
$ws.Cells.Item(9,15).Value = 0.653207634367435
$ws.Cells.Item(9,16).Value = 0.653207634367435
$ws.Cells.Item(9,17).Value = 73.69941585428268
$ws.Cells.Item(9,18).Value = 663.294742688544
$ws.Cells.Item(9,19).Value = 0.3515075705960561
$ws.Cells.Item(9,20).Value = 0.3515075705960561

$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Clcf1"
$ws.Cells.Item(10,3).Value = "Crlf1"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 4.172749333333334
$ws.Cells.Item(10,8).Value = 12.518248
$ws.Cells.Item(10,9).Value = 0.5381253250912406
$ws.Cells.Item(10,10).Value = 0.5381253250912404
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 9.293991
$ws.Cells.Item(10,14).Value = 27.881973
$ws.Cells.Item(10,15).Value = 0.3437254983469798
$ws.Cells.Item(10,16).Value = 0.3437254983469798
$ws.Cells.Item(10,17).Value = 38.78149474925601
$ws.Cells.Item(10,18).Value = 349.0334527433041
$ws.Cells.Item(10,19).Value = 0.1849673955401172
$ws.Cells.Item(10,20).Value = 0.1849673955401171
